$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows (bottom-up so row indices don't shift during deletion)
# Row 6  -> 004459875 / Helvecio / 30417.79
# Row 9  -> 008328804 / Sonia / 6000
# Row 11 -> 001761119 / Bluemetrix / 3602.36
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(6).Delete()
